$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels: "<name>_old" -> "<name>_FV2410", "<name>_new" -> "<name>_FV2504" ---
# Columns A-J hold the "_old" (FV2410) variant, columns L-U hold the "_new" (FV2504) variant;
# column K ("diff") is untouched.
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$fv2504Headers = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $fv2410Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $fv2410Headers[$i]
}

for ($i = 0; $i -lt $fv2504Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = $fv2504Headers[$i]
}

# --- 2. Turn the A1:U62 range into a proper Excel Table named "Table1" ---
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U62"), 0, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row and select the first cell below it (bottom-left pane) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
